# Add a "cliente_id" column (G) to the CUENTA_FINTECH sheet, for the
# cuenta-cliente relationship, and leave the next empty cell (I6)
# selected/formatted as on the author's machine.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Range("G1").Value = "cliente_id"

# Sequential cliente_id values for each existing row
$ws.Range("G2").Value = 1
$ws.Range("G3").Value = 2
$ws.Range("G4").Value = 3
$ws.Range("G5").Value = 4
$ws.Range("G6").Value = 5

# Stray formatted (underlined) but empty cell left by the author at I6
$ws.Range("I6").Font.Underline = $true

# Leave I6 as the active selection, matching the saved view state
[void]$ws.Range("I6").Select()

# Page setup tweak captured in the saved file
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

Write-Host "Edit applied"
